{"js": "// Update the 25 two-digit division problems in the single 20x5 table.\n// Content lives only in rows 0, 4, 8, 12, 16 (every 4th row); the other\n// rows are intentionally blank (left for pupils to write answers).\n// We target each cell by (row, col) position rather than by text search,\n// because a couple of the new values collide with OTHER cells' old values\n// (e.g. \"31\u00f74=\" -> \"29\u00f75=\" while a separate cell already holds \"29\u00f75=\" and\n// is itself being replaced) \u2014 a blind find/replace-all could clobber the\n// freshly written value. Position-based addressing sidesteps that.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"24\u00f73=\", newText: \"11\u00f73=\" },\n  { row: 0, col: 1, oldText: \"41\u00f74=\", newText: \"52\u00f76=\" },\n  { row: 0, col: 2, oldText: \"36\u00f73=\", newText: \"44\u00f77=\" },\n  { row: 0, col: 3, oldText: \"26\u00f76=\", newText: \"19\u00f72=\" },\n  { row: 0, col: 4, oldText: \"15\u00f74=\", newText: \"89\u00f79=\" },\n  { row: 4, col: 0, oldText: \"16\u00f79=\", newText: \"61\u00f74=\" },\n  { row: 4, col: 1, oldText: \"19\u00f79=\", newText: \"66\u00f77=\" },\n  { row: 4, col: 2, oldText: \"64\u00f74=\", newText: \"24\u00f76=\" },\n  { row: 4, col: 3, oldText: \"65\u00f79=\", newText: \"63\u00f77=\" },\n  { row: 4, col: 4, oldText: \"44\u00f73=\", newText: \"92\u00f73=\" },\n  { row: 8, col: 0, oldText: \"30\u00f76=\", newText: \"20\u00f72=\" },\n  { row: 8, col: 1, oldText: \"38\u00f73=\", newText: \"43\u00f74=\" },\n  { row: 8, col: 2, oldText: \"99\u00f75=\", newText: \"83\u00f77=\" },\n  { row: 8, col: 3, oldText: \"86\u00f72=\", newText: \"27\u00f79=\" },\n  { row: 8, col: 4, oldText: \"89\u00f77=\", newText: \"67\u00f76=\" },\n  { row: 12, col: 0, oldText: \"75\u00f74=\", newText: \"18\u00f77=\" },\n  { row: 12, col: 1, oldText: \"31\u00f74=\", newText: \"29\u00f75=\" },\n  { row: 12, col: 2, oldText: \"17\u00f76=\", newText: \"87\u00f78=\" },\n  { row: 12, col: 3, oldText: \"11\u00f79=\", newText: \"70\u00f74=\" },\n  { row: 12, col: 4, oldText: \"36\u00f74=\", newText: \"96\u00f73=\" },\n  { row: 16, col: 0, oldText: \"70\u00f73=\", newText: \"78\u00f74=\" },\n  { row: 16, col: 1, oldText: \"71\u00f72=\", newText: \"33\u00f73=\" },\n  { row: 16, col: 2, oldText: \"29\u00f75=\", newText: \"81\u00f73=\" },\n  { row: 16, col: 3, oldText: \"80\u00f79=\", newText: \"26\u00f73=\" },\n  { row: 16, col: 4, oldText: \"23\u00f74=\", newText: \"53\u00f78=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load every target cell's current text up front so we can sanity-check\n// before mutating (and fail loudly if the document doesn't match what we\n// expect instead of silently writing to the wrong cell).\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const { oldText, newText } = replacements[i];\n  const cell = cells[i];\n  const current = (cell.value || \"\").trim();\n  if (current !== oldText) {\n    throw new Error(\n      `Unexpected cell text at row ${replacements[i].row}, col ${replacements[i].col}: ` +\n        `expected \"${oldText}\" but found \"${current}\"`\n    );\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 two-digit division problems in the single 20x5 table.\n# Content lives only in rows 1, 5, 9, 13, 17 (1-based; every 4th row); the\n# other rows are intentionally blank (left for pupils to write answers).\n#\n# We address each cell by (row, col) position rather than by\n# Find/Replace-All on the whole document, because a couple of the NEW\n# values collide with OTHER cells' OLD values (e.g. \"31\u00f74=\" becomes\n# \"29\u00f75=\" while a different, still-unprocessed cell already holds\n# \"29\u00f75=\" as its old value) - a blind global replace-all could clobber a\n# value we just wrote. Position-based addressing sidesteps that entirely.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"24\u00f73=\"; NewText = \"11\u00f73=\" }\n    @{ Row = 1; Col = 2; OldText = \"41\u00f74=\"; NewText = \"52\u00f76=\" }\n    @{ Row = 1; Col = 3; OldText = \"36\u00f73=\"; NewText = \"44\u00f77=\" }\n    @{ Row = 1; Col = 4; OldText = \"26\u00f76=\"; NewText = \"19\u00f72=\" }\n    @{ Row = 1; Col = 5; OldText = \"15\u00f74=\"; NewText = \"89\u00f79=\" }\n    @{ Row = 5; Col = 1; OldText = \"16\u00f79=\"; NewText = \"61\u00f74=\" }\n    @{ Row = 5; Col = 2; OldText = \"19\u00f79=\"; NewText = \"66\u00f77=\" }\n    @{ Row = 5; Col = 3; OldText = \"64\u00f74=\"; NewText = \"24\u00f76=\" }\n    @{ Row = 5; Col = 4; OldText = \"65\u00f79=\"; NewText = \"63\u00f77=\" }\n    @{ Row = 5; Col = 5; OldText = \"44\u00f73=\"; NewText = \"92\u00f73=\" }\n    @{ Row = 9; Col = 1; OldText = \"30\u00f76=\"; NewText = \"20\u00f72=\" }\n    @{ Row = 9; Col = 2; OldText = \"38\u00f73=\"; NewText = \"43\u00f74=\" }\n    @{ Row = 9; Col = 3; OldText = \"99\u00f75=\"; NewText = \"83\u00f77=\" }\n    @{ Row = 9; Col = 4; OldText = \"86\u00f72=\"; NewText = \"27\u00f79=\" }\n    @{ Row = 9; Col = 5; OldText = \"89\u00f77=\"; NewText = \"67\u00f76=\" }\n    @{ Row = 13; Col = 1; OldText = \"75\u00f74=\"; NewText = \"18\u00f77=\" }\n    @{ Row = 13; Col = 2; OldText = \"31\u00f74=\"; NewText = \"29\u00f75=\" }\n    @{ Row = 13; Col = 3; OldText = \"17\u00f76=\"; NewText = \"87\u00f78=\" }\n    @{ Row = 13; Col = 4; OldText = \"11\u00f79=\"; NewText = \"70\u00f74=\" }\n    @{ Row = 13; Col = 5; OldText = \"36\u00f74=\"; NewText = \"96\u00f73=\" }\n    @{ Row = 17; Col = 1; OldText = \"70\u00f73=\"; NewText = \"78\u00f74=\" }\n    @{ Row = 17; Col = 2; OldText = \"71\u00f72=\"; NewText = \"33\u00f73=\" }\n    @{ Row = 17; Col = 3; OldText = \"29\u00f75=\"; NewText = \"81\u00f73=\" }\n    @{ Row = 17; Col = 4; OldText = \"80\u00f79=\"; NewText = \"26\u00f73=\" }\n    @{ Row = 17; Col = 5; OldText = \"23\u00f74=\"; NewText = \"53\u00f78=\" }\n)\n\nforeach ($rep in $replacements) {\n    $cell = $t.Cell($rep.Row, $rep.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $rep.OldText) {\n        throw \"Unexpected cell text at row $($rep.Row), col $($rep.Col): expected '$($rep.OldText)' but found '$current'\"\n    }\n    $cell.Range.Text = $rep.NewText\n}\n"}
